$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("variables")

# Update proxy_decision column (AC) for rows 11-13: "exclude" -> "exclude_NA"
$ws.Range("AC11:AC13").Value = "exclude_NA"

# Update variable_note column (D) for rows 5-10: "mean and variance" -> "score"
$ws.Range("D5:D10").Value = "score"

# Update the selected cell / view state to match the saved workbook view
$ws.Range("D11").Select()
